$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column O (year 2021) -----------------------------------------------
# Row 2 (bottom border strip, no value) - copy formatting from N2
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)

# Row 3 (year header)
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = 2021

# Row 4 (computed ratio, formula)
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Formula = "=O5/O6*1000"

# Row 5 (disposed waste, thousand tons) - General number format (like B5/C5)
$ws.Range("B5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 1229.5999999999999

# Row 6 (population, thousand people)
$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 6436.9

$excel.CutCopyMode = 0

# --- Recalculate so the O4 formula result is cached ----------------------
$excel.Calculate()

# --- Selection update ------------------------------------------------------
$ws.Range("P16").Select()
